$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to Text format so numeric-looking values
# (e.g. "1.003", "3.800") are preserved exactly as text, matching the source data.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "29.216.78"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3
$ws.Range("D3").Value = "1.826.09"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.39%  "

# Row 5
$ws.Range("D5").Value = "234.76"
$ws.Range("E5").Value = "  -0.67%  "

# Row 6
$ws.Range("D6").Value = "0.5995"
$ws.Range("E6").Value = "  -0.83%  "

# Row 7
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.42%  "

# Row 8
$ws.Range("D8").Value = "0.06915"
$ws.Range("E8").Value = "  -2.79%  "

# Row 9
$ws.Range("D9").Value = "0.2757"
$ws.Range("E9").Value = "  -2.23%  "

# Row 10
$ws.Range("D10").Value = "23.37"
$ws.Range("E10").Value = "  -2.79%  "

# Row 11
$ws.Range("D11").Value = "0.07609"
$ws.Range("E11").Value = "  -0.54%  "

# Row 12
$ws.Range("D12").Value = "1.831.72"
$ws.Range("E12").Value = "  +0.40%  "

# Row 13
$ws.Range("D13").Value = "4.716"
$ws.Range("E13").Value = "  -1.23%  "

# Row 14
$ws.Range("D14").Value = "0.6236"
$ws.Range("E14").Value = "  -2.63%  "

# Row 15
$ws.Range("D15").Value = "0.000009728"
$ws.Range("E15").Value = "  -2.42%  "

# Row 16
$ws.Range("D16").Value = "77.16"
$ws.Range("E16").Value = "  -3.01%  "

# Row 17
$ws.Range("D17").Value = "28.903.93"
$ws.Range("E17").Value = "  -1.10%  "

# Row 18
$ws.Range("D18").Value = "5.536"
$ws.Range("E18").Value = "  -7.78%  "

# Row 19
$ws.Range("D19").Value = "215.17"
$ws.Range("E19").Value = "  -7.13%  "

# Row 20
$ws.Range("E20").Value = "  +0.41%  "

# Row 21
$ws.Range("D21").Value = "11.53"
$ws.Range("E21").Value = "  -2.17%  "

# Row 22
$ws.Range("D22").Value = "6.814"
$ws.Range("E22").Value = "  -3.23%  "

# Row 23
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
$ws.Range("D24").Value = "155.87"
$ws.Range("E24").Value = "  +0.36%  "

# Row 25
$ws.Range("D25").Value = "7.942"
$ws.Range("E25").Value = "  -1.28%  "

# Row 26
$ws.Range("D26").Value = "0.1286"
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("D27").Value = "16.45"
$ws.Range("E27").Value = "  -1.31%  "

# Row 28
$ws.Range("D28").Value = "0.06498"
$ws.Range("E28").Value = "  -6.34%  "

# Row 29
$ws.Range("D29").Value = "1.428"
$ws.Range("E29").Value = "  -2.20%  "

# Row 30
$ws.Range("D30").Value = "1.441"
$ws.Range("E30").Value = "  -0.67%  "

# Row 31
$ws.Range("D31").Value = "3.800"
$ws.Range("E31").Value = "  -0.14%  "

# Row 32
$ws.Range("D32").Value = "3.774"
$ws.Range("E32").Value = "  -1.51%  "

# Row 33
$ws.Range("D33").Value = "1.087"
$ws.Range("E33").Value = "  -4.42%  "

# Row 34
$ws.Range("D34").Value = "1.715"
$ws.Range("E34").Value = "  -0.37%  "

# Row 35
$ws.Range("D35").Value = "0.6413"
$ws.Range("E35").Value = "  -3.37%  "

# Row 36
$ws.Range("D36").Value = "2.544"
$ws.Range("E36").Value = "  +0.51%  "

# Row 37
$ws.Range("D37").Value = "2.765"
$ws.Range("E37").Value = "  +0.42%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01747"
$ws.Range("E38").Value = "  -0.94%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "6.576"
$ws.Range("E39").Value = "  -0.24%  "

# Row 40
$ws.Range("D40").Value = "1.133.46"
$ws.Range("E40").Value = "  -8.15%  "

# Row 41
$ws.Range("D41").Value = "0.8854"
$ws.Range("E41").Value = "  -4.83%  "

# Row 42
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.31%  "

# Row 43
$ws.Range("D43").Value = "1.983.98"
$ws.Range("E43").Value = "  -0.57%  "

# Row 44
$ws.Range("D44").Value = "100.70"
$ws.Range("E44").Value = "  +0.60%  "

# Row 45
$ws.Range("D45").Value = "61.84"
$ws.Range("E45").Value = "  -2.56%  "

# Row 46
$ws.Range("E46").Value = "  -3.26%  "

# Row 47
$ws.Range("D47").Value = "1.605"
$ws.Range("E47").Value = "  -2.24%  "

# Row 48
$ws.Range("D48").Value = "8.405"
$ws.Range("E48").Value = "  -1.07%  "

# Row 49
$ws.Range("D49").Value = "0.05505"
$ws.Range("E49").Value = "  -1.43%  "

# Row 50
$ws.Range("D50").Value = "0.4536"
$ws.Range("E50").Value = "  -0.43%  "

# Row 51
$ws.Range("D51").Value = "6.339"
$ws.Range("E51").Value = "  -3.63%  "
